$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 4 days
# (quarterly Entsoe solar-production pull re-run for a later date window).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 4
}

# Updated notified-production (MW) figures from the retrained model
# for the corresponding rows.
$bUpdates = @{
    18 = 2
    19 = 2
    20 = 2
    22 = 37
    23 = 42
    24 = 49
    25 = 63
    26 = 274
    27 = 295
    28 = 322
    29 = 357
    30 = 742
    31 = 787
    32 = 846
    33 = 903
    34 = 1320
    35 = 1383
    36 = 1430
    37 = 1466
    38 = 1762
    39 = 1791
    40 = 1820
    41 = 1857
    42 = 2015
    43 = 2035
    44 = 2051
    45 = 2062
    46 = 2130
    47 = 2141
    48 = 2146
    49 = 2145
    50 = 2124
    51 = 2120
    52 = 2114
    53 = 2106
    54 = 2021
    55 = 2000
    56 = 1988
    57 = 1972
    58 = 1826
    59 = 1807
    60 = 1781
    61 = 1762
    62 = 1550
    63 = 1515
    64 = 1486
    65 = 1455
    66 = 1166
    67 = 1129
    68 = 1098
    69 = 1066
    70 = 656
    71 = 625
    72 = 594
    73 = 566
    74 = 246
    75 = 221
    76 = 201
    77 = 183
    78 = 46
    79 = 36
    80 = 30
    81 = 28
    82 = 7
    83 = 7
    84 = 7
    85 = 7
    86 = 1
    87 = 1
    88 = 1
    89 = 1
}
foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

Write-Output "done"
